$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (updated + re-sorted by total_registros descending) data for A2:B23
$data = @(
    @("FARFAN MONTOYA ROSSANA ISABEL", 71),
    @("DAVILA CORDOVA MARIBEL", 66),
    @("ANGIE BELÉN RODRÍGUEZ ZAVALA", 60),
    @("SANCHEZ ULLOA CESAR AUGUSTO", 60),
    @("RAMOS RAMOS HANDY JAIR", 59),
    @("GARCIA GUTIERREZ LUIS ARTURO", 57),
    @("CARBAJAL RAMOS JESUS MARINA", 56),
    @("FIORELA KEILY GUTIERREZ CRUZ", 56),
    @("GONZALES VALLE SEBASTIAN", 55),
    @("DELGADO DELGADO RONI", 54),
    @("VERDE LIZARRAGA DEYSI EUFEMIA", 53),
    @("OLIVA ALVA GOSSELYN NASSIRA", 52),
    @("BAZAN TEJADA JOSE VICENTE", 52),
    @("DANY DARWIN VILLACORTA SAAVEDRA", 51),
    @("ARANEDA LOPEZ MARCO VIERI", 50),
    @("CASTILLO QUEZADA DIEGO ALONSO", 49),
    @("JOSSY IVANA SUÁREZ ZAVALETA", 49),
    @("GUERRA CALDERON ESTHEFANY NICOLLE", 47),
    @("JAVE CHAVEZ ANGHELO MARTIN", 43),
    @("CARDENAS CAMPOJO MARY PAULA", 43),
    @("ALVITES CAMPOS SERGIO MARTIN", 21),
    @("ROCHA SIPIRAN JHORDAN ENRIQUE", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
